# Generate Report for Handoff
# The c8686c97-d533-4f47-adbc-7be43f0dea3d.md file has finished being
# handed back and is now ready to be handed off again, so update its
# status to "Ready for handoff" on every sheet, and stamp the new
# handoff datetime for each locale on its respective status sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "2016-02-24 06:51:18"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("D3").Value = "2016-02-24 06:51:34"
